$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Create")
$ws.Range("C4").Value = "f"
$ws.Range("C5").Value = "#f"
$ws.Range("C6").Value = "f"
